# Generate Report for Handback
# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# values on the zh-cn and de-de worksheets to reflect the latest report run.

$wb = $excel.ActiveWorkbook

# zh-cn sheet: row 2 corresponds to the
# ef437b09-99f8-4ab4-8898-f97ebbf73d97.dff64b847da2c97857111235986827dfc6f9982e.zh-cn.xlf record
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-12 08:46:57"
$wsZhCn.Range("H2").Value = "2016-03-12 08:47:14"

# de-de sheet: row 2 corresponds to the
# ef437b09-99f8-4ab4-8898-f97ebbf73d97.dff64b847da2c97857111235986827dfc6f9982e.de-de.xlf record
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-12 08:47:00"
$wsDeDe.Range("H2").Value = "2016-03-12 08:47:19"
